$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$passColor = 5296274   # RGB(146, 208, 80) -> matches existing fgColor FF92D050 (Pass)
$failColor = 255       # RGB(255, 0, 0)    -> matches existing fgColor FFFF0000 (Fail)

$passRows = @(4, 5, 6, 7, 9, 10)
foreach ($r in $passRows) {
    $cell = $ws.Range("G$r")
    $cell.Value = "Pass"
    $cell.Interior.Color = $passColor
}

$failCell = $ws.Range("G8")
$failCell.Value = "Fail"
$failCell.Interior.Color = $failColor

$ws.Range("G10").Select() | Out-Null
